$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing "Test" label and comment text (currently in K2 / K5)
# before they get overwritten by the new LTSD parameters table, so they
# can be relocated to column O.
$oldK2 = $ws.Range("K2").Value2
$oldK5 = $ws.Range("K5").Value2

# --- New LTSD Parameters table (columns K:N, rows 2:5) ---------------

# Header
$ws.Range("K2").Value = "LTSD Parameters"

# Group headers
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

# Sub-headers ("Win" is entered for both groups before the remaining
# sub-headers, matching the original authoring order)
$ws.Range("L4").Value = "Win"
$ws.Range("N4").Value = "Win"
$ws.Range("K4").Value = "Threshols"
$ws.Range("M4").Value = "Threshold"

# Values (kept as text, not numbers, to match the source data)
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "200.0"
$ws.Range("N5").Style = "Normal"

$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = "4.3"
$ws.Range("K5").Style = "Normal"

$ws.Range("L5").NumberFormat = "@"
$ws.Range("L5").Value = "100.0"
$ws.Range("L5").Style = "Normal"

$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "5.0"
$ws.Range("M5").Style = "Normal"

# --- Relocate the pre-existing "Test" label / comment to column O ----
$ws.Range("O2").Value = $oldK2
$ws.Range("O5").Value = $oldK5

# --- Update the selected cell shown when the workbook is reopened ----
$ws.Range("M5").Select() | Out-Null
